$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 405.125
$ws.Range("I2").Value = 405.125
$ws.Range("K2").Value = 405.125
$ws.Range("M2").Value = -292.125
$ws.Range("H33").Value = 183.6
$ws.Range("I33").Value = 150.83333
$ws.Range("K33").Value = 150.83333
$ws.Range("M33").Value = 78.16667000000001
$ws.Range("H40").Value = 3339.4
$ws.Range("I40").Value = 3199
$ws.Range("K40").Value = 3199
$ws.Range("M40").Value = -3024
$ws.Range("H41").Value = 674.75
$ws.Range("I41").Value = 499.66666
$ws.Range("K41").Value = 499.66666
$ws.Range("M41").Value = -59.66665999999998
$ws.Range("H64").Value = 19510
$ws.Range("J64").Value = 19314
$ws.Range("L64").Value = 19314
$ws.Range("N64").Value = -19810
$ws.Range("H67").Value = 19510
$ws.Range("J67").Value = 19314
$ws.Range("L67").Value = 19314
$ws.Range("N67").Value = -21030
$ws.Range("H86").Value = 5510.727
$ws.Range("I86").Value = 5381.1665
$ws.Range("J86").Value = 5666.2
$ws.Range("K86").Value = 5381.1665
$ws.Range("L86").Value = 5666.2
$ws.Range("M86").Value = -4258.1665
$ws.Range("N86").Value = -7912.2
$ws.Range("H89").Value = 5510.727
$ws.Range("I89").Value = 5381.1665
$ws.Range("J89").Value = 5666.2
$ws.Range("K89").Value = 26905.8325
$ws.Range("L89").Value = 28331
$ws.Range("M89").Value = -21289.8325
$ws.Range("N89").Value = -39563
$ws.Range("H96").Value = 1515.5
$ws.Range("I96").Value = 531.25
$ws.Range("J96").Value = 2499.75
$ws.Range("K96").Value = 1593.75
$ws.Range("L96").Value = 7499.25
$ws.Range("M96").Value = -220.75
$ws.Range("N96").Value = -10245.25
$ws.Range("H97").Value = 5599.4
$ws.Range("J97").Value = 5599.4
$ws.Range("L97").Value = 16798.2
$ws.Range("N97").Value = -17790.2
$ws.Range("H100").Value = 2473.4
$ws.Range("I100").Value = 2473.4
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2473.4
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1932.4
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 10983
$ws.Range("I113").Value = 11263.143
$ws.Range("J113").Value = 10002.5
$ws.Range("K113").Value = 11263.143
$ws.Range("L113").Value = 10002.5
$ws.Range("M113").Value = -8009.143
$ws.Range("N113").Value = -16510.5
$ws.Range("H135").Value = 1192.6666
$ws.Range("I135").Value = 431.3
$ws.Range("K135").Value = 3881.7
$ws.Range("M135").Value = -1346.7
$ws.Range("H138").Value = 3060
$ws.Range("I138").Value = 2482.2222
$ws.Range("J138").Value = 4100
$ws.Range("K138").Value = 7446.6666
$ws.Range("L138").Value = 12300
$ws.Range("M138").Value = -2306.6666
$ws.Range("N138").Value = -22580
$ws.Range("H141").Value = 922.75
$ws.Range("I141").Value = 934
$ws.Range("J141").Value = 799
$ws.Range("K141").Value = 2802
$ws.Range("L141").Value = 2397
$ws.Range("M141").Value = 2378
$ws.Range("N141").Value = -12757

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 11073
$ws.Range("I28").Value = 11073
$ws.Range("K28").Value = 11073
$ws.Range("M28").Value = -10881
$ws.Range("H32").Value = 2698.4783
$ws.Range("I32").Value = 2316.5264
$ws.Range("J32").Value = 4512.75
$ws.Range("K32").Value = 2316.5264
$ws.Range("L32").Value = 4512.75
$ws.Range("M32").Value = -2029.5264
$ws.Range("N32").Value = -5086.75
$ws.Range("H45").Value = 3592.7856
$ws.Range("I45").Value = 1708.8
$ws.Range("J45").Value = 4639.4443
$ws.Range("K45").Value = 1708.8
$ws.Range("L45").Value = 4639.4443
$ws.Range("M45").Value = -1331.8
$ws.Range("N45").Value = -5393.4443
$ws.Range("H63").Value = 2250
$ws.Range("I63").Value = 2250
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2250
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1564
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2250
$ws.Range("I66").Value = 2250
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 11250
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -7818
$ws.Range("N66").ClearContents()
$ws.Range("H99").Value = 11073
$ws.Range("I99").Value = 11073
$ws.Range("K99").Value = 11073
$ws.Range("M99").Value = -8078
$ws.Range("H110").Value = 1218
$ws.Range("I110").Value = 1218
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1218
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 827
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 1313.75
$ws.Range("I122").Value = 1079.4
$ws.Range("J122").Value = 1704.3334
$ws.Range("K122").Value = 3238.2
$ws.Range("L122").Value = 5113.0002
$ws.Range("M122").Value = -788.2000000000003
$ws.Range("N122").Value = -10013.0002
$ws.Range("H132").Value = 1850
$ws.Range("I132").Value = 1850
$ws.Range("K132").Value = 5550
$ws.Range("M132").Value = -3020
$ws.Range("H133").Value = 49999
$ws.Range("J133").Value = 49999
$ws.Range("L133").Value = 49999
$ws.Range("N133").Value = -55059

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3131.7334
$ws.Range("I86").Value = 2817.375
$ws.Range("K86").Value = 2817.375
$ws.Range("M86").Value = -1694.375
$ws.Range("H89").Value = 3131.7334
$ws.Range("I89").Value = 2817.375
$ws.Range("K89").Value = 14086.875
$ws.Range("M89").Value = -8470.875
$ws.Range("H94").Value = 2059.9
$ws.Range("I94").Value = 799.8570999999999
$ws.Range("K94").Value = 799.8570999999999
$ws.Range("M94").Value = -348.8570999999999
$ws.Range("H99").Value = 1915.5555
$ws.Range("I99").Value = 1672.0667
$ws.Range("K99").Value = 1672.0667
$ws.Range("M99").Value = -174.0667000000001
$ws.Range("H134").Value = 1152.381
$ws.Range("I134").Value = 1152.381
$ws.Range("K134").Value = 3457.143
$ws.Range("M134").Value = -922.143
$ws.Range("H137").Value = 89998
$ws.Range("J137").Value = 89998
$ws.Range("L137").Value = 89998
$ws.Range("N137").Value = -100198

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 466.33334
$ws.Range("I16").Value = 450
$ws.Range("J16").Value = 499
$ws.Range("K16").Value = 450
$ws.Range("L16").Value = 499
$ws.Range("M16").Value = -163
$ws.Range("N16").Value = -1073
$ws.Range("H31").Value = 4996.75
$ws.Range("J31").Value = 4997.6665
$ws.Range("L31").Value = 4997.6665
$ws.Range("N31").Value = -5587.6665
$ws.Range("H34").Value = 4996.75
$ws.Range("J34").Value = 4997.6665
$ws.Range("L34").Value = 4997.6665
$ws.Range("N34").Value = -5401.6665
$ws.Range("H58").Value = 1037
$ws.Range("I58").Value = 1037
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1037
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -834
$ws.Range("N58").ClearContents()
$ws.Range("H113").Value = 466.33334
$ws.Range("I113").Value = 450
$ws.Range("J113").Value = 499
$ws.Range("K113").Value = 450
$ws.Range("L113").Value = 499
$ws.Range("M113").Value = 1720
$ws.Range("N113").Value = -4839
$ws.Range("H122").Value = 2624.6
$ws.Range("I122").Value = 2624.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7873.799999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5423.799999999999
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2812.8462
$ws.Range("I132").Value = 2812.8462
$ws.Range("K132").Value = 8438.5386
$ws.Range("M132").Value = -5908.5386
$ws.Range("H136").Value = 1037
$ws.Range("I136").Value = 1037
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3111
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -561
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 340.4
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 425.25
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1275.75
$ws.Range("M12").Value = 170
$ws.Range("N12").Value = -1621.75
$ws.Range("H122").Value = 387
$ws.Range("J122").Value = 400
$ws.Range("L122").Value = 3600
$ws.Range("N122").Value = -8500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 928562.5
$ws.Range("I11").Value = 403916.66
$ws.Range("J11").Value = 2502500
$ws.Range("K11").Value = 403916.66
$ws.Range("L11").Value = 2502500
$ws.Range("M11").Value = -403777.66
$ws.Range("N11").Value = -2502778
$ws.Range("H12").Value = 3111555.5
$ws.Range("I12").Value = 3111555.5
$ws.Range("K12").Value = 3111555.5
$ws.Range("M12").Value = -3111415.5
$ws.Range("H18").Value = 26999.5
$ws.Range("J18").Value = 49999
$ws.Range("L18").Value = 49999
$ws.Range("N18").Value = -50585
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H70").Value = 10999.75
$ws.Range("I70").Value = 10999.75
$ws.Range("K70").Value = 10999.75
$ws.Range("M70").Value = -10729.75
$ws.Range("H73").Value = 10999.75
$ws.Range("I73").Value = 10999.75
$ws.Range("K73").Value = 10999.75
$ws.Range("M73").Value = -10063.75
$ws.Range("H97").Value = 3750
$ws.Range("I97").Value = 3000
$ws.Range("K97").Value = 3000
$ws.Range("M97").Value = -2504
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 1754.5
$ws.Range("I122").Value = 1617.3334
$ws.Range("J122").Value = 2001.4
$ws.Range("K122").Value = 4852.0002
$ws.Range("L122").Value = 6004.200000000001
$ws.Range("M122").Value = -2402.0002
$ws.Range("N122").Value = -10904.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 1336.6666
$ws.Range("H61").Value = 1695
$ws.Range("I61").Value = 1695
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1695
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1493
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 5000
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -6498
$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 25000
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -32488
$ws.Range("H99").Value = 20498.334
$ws.Range("I99").Value = 20498.334
$ws.Range("K99").Value = 20498.334
$ws.Range("M99").Value = -17503.334
$ws.Range("H100").Value = 2083.3333
$ws.Range("I100").Value = 2083.3333
$ws.Range("K100").Value = 2083.3333
$ws.Range("M100").Value = -1542.3333
$ws.Range("H113").Value = 1695
$ws.Range("I113").Value = 1695
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1695
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 475
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1209.4
$ws.Range("I96").Value = 1099.1428
$ws.Range("J96").Value = 1466.6666
$ws.Range("K96").Value = 1099.1428
$ws.Range("L96").Value = 1466.6666
$ws.Range("M96").Value = 273.8571999999999
$ws.Range("N96").Value = -4212.6666
$ws.Range("H100").Value = 3873010.8
$ws.Range("I100").Value = 4647319.5
$ws.Range("J100").Value = 1466.6666
$ws.Range("K100").Value = 9294639
$ws.Range("L100").Value = 2933.3332
$ws.Range("M100").Value = -9294098
$ws.Range("N100").Value = -4015.3332
$ws.Range("H132").Value = 1908.7778
$ws.Range("I132").Value = 1647.375
$ws.Range("K132").Value = 4942.125
$ws.Range("M132").Value = -2412.125
